# Update gh-pages to output generated at 456a3b4
# Bumps a handful of "想去人数" (interest count) values across sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 2277
$wsExhibit.Range("F8").Value = 1711
$wsExhibit.Range("F9").Value = 2886
$wsExhibit.Range("F10").Value = 162
$wsExhibit.Range("F11").Value = 4290
$wsExhibit.Range("F20").Value = 97
$wsExhibit.Range("F21").Value = 294
$wsExhibit.Range("F22").Value = 4137
$wsExhibit.Range("F24").Value = 3621
$wsExhibit.Range("F25").Value = 1126
$wsExhibit.Range("F27").Value = 534
$wsExhibit.Range("F31").Value = 511
$wsExhibit.Range("F32").Value = 450

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 25

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 2277
$wsAll.Range("F10").Value = 1711
$wsAll.Range("F12").Value = 2886
$wsAll.Range("F13").Value = 162
$wsAll.Range("F14").Value = 4290
$wsAll.Range("F22").Value = 25
$wsAll.Range("F24").Value = 97
$wsAll.Range("F25").Value = 294
$wsAll.Range("F26").Value = 4137
$wsAll.Range("F28").Value = 3621
$wsAll.Range("F29").Value = 1126
$wsAll.Range("F31").Value = 534
$wsAll.Range("F35").Value = 511
$wsAll.Range("F36").Value = 450
